$d = $word.ActiveDocument

# NOTE: Range.Find.Execute() collapses the range to the found/replaced text,
# so each call below re-fetches a fresh paragraph Range (scoped search) to
# keep later replacements within the same paragraph working correctly, and
# to avoid touching look-alike text elsewhere in the document.

function Replace-InParagraph($index, $old, $new) {
    $r = $d.Paragraphs.Item($index).Range
    $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Paragraph 1: "English / Portuguese / French / Thai / Vietnamese / Spanish" ---
Replace-InParagraph 1 "English" "الإنجليزية"
Replace-InParagraph 1 " / Portuguese / French / Thai / Vietnamese / Spanish" " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية"

# --- Paragraph 3: "English" (standalone) ---
Replace-InParagraph 3 "English" "الإنجليزية"

# --- Paragraph 5: "Brief:" ---
Replace-InParagraph 5 "Brief" "المضمون"

# --- Paragraph 6: brief description ---
Replace-InParagraph 6 "An email sent to partners who have attended the event. This email will include a photo gallery It will be sent via customer.io" "تم إرسال بريد إلكتروني إلى الشركاء الذين حضروا الحدث. سيتضمن هذا البريد الإلكتروني معرض صور سيتم إرساله عبر customer.io"

# --- Paragraph 8: "Target audience:" ---
Replace-InParagraph 8 "Target audience" "الجمهور المستهدف"

# --- Paragraph 9: "Event attendees" ---
Replace-InParagraph 9 "Event attendees" "الحاضرون في الحدث"

# --- Paragraph 12: "Subject: Thank you for coming to [EVENT NAME]! " ---
Replace-InParagraph 12 "Subject: " "الموضوع: "
Replace-InParagraph 12 "Thank you for coming to " "شكرًا لقدومك إلى "
Replace-InParagraph 12 "[EVENT NAME]" "[اسم الحدث]"

# --- Paragraph 13: "You made our event a success! 🎉" ---
Replace-InParagraph 13 "You made our event a success! 🎉" "لقد ساهمت في نجاح هذا الحدث! 🎉"

# --- Paragraph 15: "Hi [PARTNER NAME], " ---
Replace-InParagraph 15 "Hi " "مرحبًا "
Replace-InParagraph 15 "[PARTNER NAME]" "[اسم الشريك]"
Replace-InParagraph 15 ", " "، "

# --- Paragraph 17: "Thank you for attending [EVENT NAME] in [CITY], [COUNTRY]. We hope you had a great time, and it was a pleasure getting to know you!" ---
Replace-InParagraph 17 "Thank you for attending " "شكرًا لحضورك "
Replace-InParagraph 17 "[EVENT NAME]" "[اسم الحدث]"
Replace-InParagraph 17 " in " " في "
Replace-InParagraph 17 "[CITY]" "[المدينة]"

# Only the comma right after [CITY] (before [COUNTRY]) must change; the later
# comma in "a great time, and it was a pleasure..." must stay untouched, so
# scope the Find to a fresh sub-range ending right after "[COUNTRY]".
$p17 = $d.Paragraphs.Item(17).Range
$fullText = $p17.Text
$countryIdx = $fullText.IndexOf("[COUNTRY]")
$subStart = $p17.Start
$subEnd = $subStart + $countryIdx + 9
$citySub = $d.Range($subStart, $subEnd)
$citySub.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, "، ", 2) | Out-Null

$p17b = $d.Paragraphs.Item(17).Range
$fullText2 = $p17b.Text
$countryIdx2 = $fullText2.IndexOf("[COUNTRY]")
$subStart2 = $p17b.Start
$subEnd2 = $subStart2 + $countryIdx2 + 9
$citySub2 = $d.Range($subStart2, $subEnd2)
$citySub2.Find.Execute("[COUNTRY]", $true, $false, $false, $false, $false, $true, 1, $false, "[البلد]", 2) | Out-Null
